# Update cryptocurrency price/volume figures per the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "332.03"
Set-TextValue $ws.Range("E2") "0.99%"

# Row 3
Set-TextValue $ws.Range("D3") "41.56"
Set-TextValue $ws.Range("E3") "2.48%"

# Row 4
Set-TextValue $ws.Range("D4") "5.700"
Set-TextValue $ws.Range("E4") "-3.53%"

# Row 5
Set-TextValue $ws.Range("D5") "0.08107"
Set-TextValue $ws.Range("E5") "-0.36%"

# Row 6
Set-TextValue $ws.Range("D6") "2.060"
Set-TextValue $ws.Range("E6") "5.32%"

# Row 7
Set-TextValue $ws.Range("D7") "8.734"
Set-TextValue $ws.Range("E7") "-0.07%"

# Row 8
Set-TextValue $ws.Range("D8") "4.527"
Set-TextValue $ws.Range("E8") "-1.31%"

# Row 9
Set-TextValue $ws.Range("D9") "2.973"
Set-TextValue $ws.Range("E9") "1.02%"

# Row 10
Set-TextValue $ws.Range("D10") "0.9253"
Set-TextValue $ws.Range("E10") "-2.01%"

# Row 11
Set-TextValue $ws.Range("D11") "0.1265"
Set-TextValue $ws.Range("E11") "-3.27%"

# Row 12
Set-TextValue $ws.Range("D12") "0.1963"
Set-TextValue $ws.Range("E12") "-1.65%"

# Row 13
Set-TextValue $ws.Range("D13") "8.810"
Set-TextValue $ws.Range("E13") "13.81%"

# Row 14
Set-TextValue $ws.Range("D14") "0.09186"
Set-TextValue $ws.Range("E14") "-1.05%"

# Row 15
Set-TextValue $ws.Range("D15") "0.03699"
Set-TextValue $ws.Range("E15") "8.35%"

# Row 16
Set-TextValue $ws.Range("D16") "0.1051"
Set-TextValue $ws.Range("E16") "9.37%"

# Row 17
Set-TextValue $ws.Range("D17") "0.001308"

# Row 18
Set-TextValue $ws.Range("D18") "0.006153"
Set-TextValue $ws.Range("E18") "0.87%"

# Row 19
Set-TextValue $ws.Range("D19") "3.380"
Set-TextValue $ws.Range("E19") "0.17%"

# Row 20
Set-TextValue $ws.Range("D20") "0.3531"
Set-TextValue $ws.Range("E20") "1.00%"

# Row 21
Set-TextValue $ws.Range("D21") "0.1417"
Set-TextValue $ws.Range("E21") "-2.09%"

# Row 22
Set-TextValue $ws.Range("D22") "0.2612"
Set-TextValue $ws.Range("E22") "6.67%"

# Row 23
Set-TextValue $ws.Range("E23") "0.03%"

# Row 24
Set-TextValue $ws.Range("D24") "0.001257"
Set-TextValue $ws.Range("E24") "0.21%"

# Row 25
Set-TextValue $ws.Range("D25") "0.004455"
Set-TextValue $ws.Range("E25") "2.29%"

# Row 26
Set-TextValue $ws.Range("D26") "0.0001240"
Set-TextValue $ws.Range("E26") "4.15%"

# Row 39
Set-TextValue $ws.Range("D39") "0.02831"
Set-TextValue $ws.Range("E39") "13.27%"

# Row 40
Set-TextValue $ws.Range("D40") "0.05506"
Set-TextValue $ws.Range("E40") "4.12%"

# Row 41
Set-TextValue $ws.Range("D41") "0.007651"
Set-TextValue $ws.Range("E41") "0.78%"

# Row 42
Set-TextValue $ws.Range("D42") "0.009918"
Set-TextValue $ws.Range("E42") "10.78%"

# Row 43
Set-TextValue $ws.Range("D43") "0.1419"
Set-TextValue $ws.Range("E43") "-0.85%"

# Row 44
Set-TextValue $ws.Range("D44") "0.002090"
Set-TextValue $ws.Range("E44") "1.13%"

# Row 45
Set-TextValue $ws.Range("D45") "0.01075"
Set-TextValue $ws.Range("E45") "13.31%"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006818"
Set-TextValue $ws.Range("E46") "-0.26%"

# Row 47
Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "0.11%"

# Row 48
Set-TextValue $ws.Range("D48") "0.002996"
Set-TextValue $ws.Range("E48") "3.33%"

# Row 49
Set-TextValue $ws.Range("D49") "0.002280"
Set-TextValue $ws.Range("E49") "26.75%"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("E50") "0.11%"

# Row 51
Set-TextValue $ws.Range("D51") "0.0002001"
Set-TextValue $ws.Range("E51") "0.11%"
